$wb = $excel.ActiveWorkbook

# Rename "Repeaters" -> "RepeatersOld" first so the second rename below
# (which reuses the name "Repeaters") doesn't collide with it.
$wsOld = $wb.Worksheets.Item("Repeaters")
$wsOld.Name = "RepeatersOld"

# "Repeaters Updated" becomes the new "Repeaters" sheet.
$wsNew = $wb.Worksheets.Item("Repeaters Updated")
$wsNew.Name = "Repeaters"

# Make the (renamed) "Repeaters" sheet the active tab, with A15 selected.
$wsNew.Activate()
$wsNew.Range("A15").Select()
